# Update the "Good Morning" greeting cell to "GIT UPDATE" on the Rules sheet,
# and leave the selection on that cell (matches the author's edit location).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")
$ws.Range("E8").Value = "GIT UPDATE"
$ws.Range("E8").Select()
